$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 451.81818
$ws.Range("J55").Value = 130
$ws.Range("L55").Value = 130
$ws.Range("N55").Value = -558

$ws.Range("H107").Value = 78125530
$ws.Range("I107").Value = 20833968
$ws.Range("J107").Value = 250000210
$ws.Range("K107").Value = 20833968
$ws.Range("L107").Value = 250000210
$ws.Range("M107").Value = -20832048
$ws.Range("N107").Value = -250004050

$ws.Range("H132").Value = 807.9231
$ws.Range("I132").Value = 807.9231
$ws.Range("K132").Value = 2423.7693
$ws.Range("M132").Value = 106.2307000000001

$ws.Range("H137").Value = 6134.9287
$ws.Range("I137").Value = 3707.5833
$ws.Range("K137").Value = 11122.7499
$ws.Range("M137").Value = -8572.749899999999

$ws.Range("H138").Value = 1153059.2
$ws.Range("I138").Value = 2001.7878
$ws.Range("J138").Value = 1856483.2
$ws.Range("K138").Value = 6005.3634
$ws.Range("L138").Value = 5569449.6
$ws.Range("M138").Value = -865.3634000000002
$ws.Range("N138").Value = -5579729.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3082923.5
$ws.Range("I32").Value = 3130313
$ws.Range("K32").Value = 3130313
$ws.Range("M32").Value = -3130026

$ws.Range("H45").Value = 5119.222
$ws.Range("I45").Value = 1863.6666
$ws.Range("J45").Value = 6747
$ws.Range("K45").Value = 1863.6666
$ws.Range("L45").Value = 6747
$ws.Range("M45").Value = -1486.6666
$ws.Range("N45").Value = -7501

$ws.Range("H61").Value = 100013140
$ws.Range("I61").Value = 5668.5
$ws.Range("K61").Value = 5668.5
$ws.Range("M61").Value = -5456.5

$ws.Range("H74").Value = 2702.6047
$ws.Range("I74").Value = 1752.258
$ws.Range("K74").Value = 1752.258
$ws.Range("M74").Value = -878.258

$ws.Range("H77").Value = 2702.6047
$ws.Range("I77").Value = 1752.258
$ws.Range("K77").Value = 8761.290000000001
$ws.Range("M77").Value = -4393.290000000001

$ws.Range("H132").Value = 5420.0264
$ws.Range("I132").Value = 1717.1428
$ws.Range("J132").Value = 9994.177
$ws.Range("K132").Value = 5151.428400000001
$ws.Range("L132").Value = 29982.531
$ws.Range("M132").Value = -2621.428400000001
$ws.Range("N132").Value = -35042.531

$ws.Range("H136").Value = 100013140
$ws.Range("I136").Value = 5668.5
$ws.Range("K136").Value = 17005.5
$ws.Range("M136").Value = -14455.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 41585.5
$ws.Range("I51").Value = 13999
$ws.Range("J51").Value = 50781
$ws.Range("K51").Value = 13999
$ws.Range("L51").Value = 50781
$ws.Range("M51").Value = -13508
$ws.Range("N51").Value = -51763

$ws.Range("H99").Value = 3249316
$ws.Range("I99").Value = 2288.889
$ws.Range("J99").Value = 9093965
$ws.Range("K99").Value = 2288.889
$ws.Range("L99").Value = 9093965
$ws.Range("M99").Value = -790.8890000000001
$ws.Range("N99").Value = -9096961

$ws.Range("H105").Value = 2918
$ws.Range("I105").Value = 2000.2222
$ws.Range("K105").Value = 2000.2222
$ws.Range("M105").Value = -253.2221999999999

$ws.Range("H107").Value = 562500000
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 8935804
$ws.Range("I134").Value = 20835502
$ws.Range("K134").Value = 62506506
$ws.Range("M134").Value = -62503971

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6265.385
$ws.Range("I16").Value = 3592.2
$ws.Range("J16").Value = 7936.125
$ws.Range("K16").Value = 3592.2
$ws.Range("L16").Value = 7936.125
$ws.Range("M16").Value = -3305.2
$ws.Range("N16").Value = -8510.125

$ws.Range("H25").Value = 782.2
$ws.Range("I25").Value = 477.75
$ws.Range("J25").Value = 2000
$ws.Range("K25").Value = 477.75
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = -303.75
$ws.Range("N25").Value = -2348

$ws.Range("H31").Value = 7065.447
$ws.Range("I31").Value = 2389.4
$ws.Range("K31").Value = 2389.4
$ws.Range("M31").Value = -2094.4

$ws.Range("H34").Value = 7065.447
$ws.Range("I34").Value = 2389.4
$ws.Range("K34").Value = 2389.4
$ws.Range("M34").Value = -2187.4

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").ClearContents()
$ws.Range("N111").Value = 0

$ws.Range("H113").Value = 6265.385
$ws.Range("I113").Value = 3592.2
$ws.Range("J113").Value = 7936.125
$ws.Range("K113").Value = 3592.2
$ws.Range("L113").Value = 7936.125
$ws.Range("M113").Value = -1422.2
$ws.Range("N113").Value = -12276.125

$ws.Range("H122").Value = 1731.5555
$ws.Range("I122").Value = 936.2857
$ws.Range("J122").Value = 2237.6365
$ws.Range("K122").Value = 2808.8571
$ws.Range("L122").Value = 6712.9095
$ws.Range("M122").Value = -358.8571000000002
$ws.Range("N122").Value = -11612.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 53277268
$ws.Range("I4").Value = 61063736
$ws.Range("K4").Value = 183191208
$ws.Range("M4").Value = -183191096

$ws.Range("H137").Value = 3825
$ws.Range("I137").Value = 3258.6667
$ws.Range("J137").Value = 4674.5
$ws.Range("K137").Value = 9776.000100000001
$ws.Range("L137").Value = 14023.5
$ws.Range("M137").Value = -4676.000100000001
$ws.Range("N137").Value = -24223.5

$ws.Range("H138").Value = 7206.5625
$ws.Range("I138").Value = 6307.5
$ws.Range("K138").Value = 18922.5
$ws.Range("M138").Value = -13782.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1465.0952
$ws.Range("I97").Value = 757.1667
$ws.Range("J97").Value = 2409
$ws.Range("K97").Value = 757.1667
$ws.Range("L97").Value = 2409
$ws.Range("M97").Value = -261.1667
$ws.Range("N97").Value = -3401

$ws.Range("H102").Value = 2806.4
$ws.Range("I102").Value = 2960.3845
$ws.Range("K102").Value = 2960.3845
$ws.Range("M102").Value = -1338.3845

$ws.Range("H132").Value = 5054.5864
$ws.Range("I132").Value = 1989.5555
$ws.Range("J132").Value = 10070.091
$ws.Range("K132").Value = 5968.666499999999
$ws.Range("L132").Value = 30210.273
$ws.Range("M132").Value = -3438.666499999999
$ws.Range("N132").Value = -35270.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3982.6191
$ws.Range("I7").Value = 2977.1875
$ws.Range("K7").Value = 2977.1875
$ws.Range("M7").Value = -2865.1875

$ws.Range("H22").Value = 1873

$ws.Range("H27").Value = 1873

$ws.Range("H40").Value = 2884.4
$ws.Range("I40").Value = 1954.129
$ws.Range("K40").Value = 1954.129
$ws.Range("M40").Value = -1818.129

$ws.Range("H46").Value = 1381416.9
$ws.Range("I46").Value = 2300265.2
$ws.Range("J46").Value = 3144.2
$ws.Range("K46").Value = 2300265.2
$ws.Range("L46").Value = 3144.2
$ws.Range("M46").Value = -2300077.2
$ws.Range("N46").Value = -3520.2

$ws.Range("H93").Value = 1585.875
$ws.Range("I93").Value = 937.4
$ws.Range("K93").Value = 937.4
$ws.Range("M93").Value = 310.6

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0

$ws.Range("H126").Value = 3982.6191
$ws.Range("I126").Value = 2977.1875
$ws.Range("K126").Value = 8931.5625
$ws.Range("M126").Value = -6461.5625

$ws.Range("H132").Value = 11117600
$ws.Range("I132").Value = 23812238
$ws.Range("K132").Value = 71436714
$ws.Range("M132").Value = -71434184

$ws.Range("H136").Value = 10037.36
$ws.Range("I136").Value = 3035.1428
$ws.Range("J136").Value = 15978.637
$ws.Range("K136").Value = 9105.428400000001
$ws.Range("L136").Value = 47935.911
$ws.Range("M136").Value = -6555.428400000001
$ws.Range("N136").Value = -53035.911

$ws.Range("H140").Value = 69407.42999999999
$ws.Range("J140").Value = 69407.42999999999
$ws.Range("L140").Value = 69407.42999999999
$ws.Range("N140").Value = -79767.42999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 28333.334
$ws.Range("I51").Value = 27500
$ws.Range("J51").Value = 30000
$ws.Range("K51").Value = 27500
$ws.Range("L51").Value = 30000
$ws.Range("M51").Value = -26990
$ws.Range("N51").Value = -31020

$ws.Range("H100").Value = 591.88
$ws.Range("I100").Value = 386.25
$ws.Range("J100").Value = 957.44446
$ws.Range("K100").Value = 772.5
$ws.Range("L100").Value = 1914.88892
$ws.Range("M100").Value = -231.5
$ws.Range("N100").Value = -2996.88892

$ws.Range("H107").Value = 11905714
$ws.Range("I107").Value = 442.10526
$ws.Range("J107").Value = 37039068
$ws.Range("K107").Value = 1326.31578
$ws.Range("L107").Value = 111117204
$ws.Range("M107").Value = 593.6842200000001
$ws.Range("N107").Value = -111121044

$ws.Range("H113").Value = 11726.4
$ws.Range("I113").Value = 22295.666
$ws.Range("K113").Value = 66886.99800000001
$ws.Range("M113").Value = -64716.99800000001

$ws.Range("H122").Value = 254103.69
$ws.Range("I122").Value = 366232.2
$ws.Range("K122").Value = 1098696.6
$ws.Range("M122").Value = -1096246.6

$ws.Range("H126").Value = 738.4
$ws.Range("I126").Value = 966.3333
$ws.Range("K126").Value = 2898.9999
$ws.Range("M126").Value = -428.9998999999998

$ws.Range("H132").Value = 10662.059
$ws.Range("I132").Value = 15256.5
$ws.Range("J132").Value = 6578.1113
$ws.Range("K132").Value = 45769.5
$ws.Range("L132").Value = 19734.3339
$ws.Range("M132").Value = -43239.5
$ws.Range("N132").Value = -24794.3339

$ws.Range("H136").Value = 48101636
$ws.Range("J136").Value = 562908.9
$ws.Range("L136").Value = 1688726.7
$ws.Range("N136").Value = -1693826.7

$ws.Range("H140").Value = 65453.5
$ws.Range("J140").Value = 69089.71000000001
$ws.Range("L140").Value = 69089.71000000001
$ws.Range("N140").Value = -79449.71000000001
